$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "A1"
$ws.Range("C2").Value = "A1"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = "M1"
$ws.Range("F2").Value = "M1"
$ws.Range("G2").Value = "M3"
$ws.Range("J2").Value = "DO"
$ws.Range("K2").Value = "M1"
$ws.Range("L2").Value = "A1"
$ws.Range("M2").Value = "A1"
$ws.Range("O2").Value = "M3"
$ws.Range("R2").Value = "M1"
$ws.Range("S2").Value = "A1"
$ws.Range("T2").Value = "DO"
$ws.Range("U2").Value = "M3"
$ws.Range("X2").Value = "M3"
$ws.Range("Y2").Value = "M1"
$ws.Range("Z2").Value = "M1"
$ws.Range("AA2").Value = "DO"
$ws.Range("B3").Value = "M1"
$ws.Range("C3").Value = "DO"
$ws.Range("D3").Value = "M1"
$ws.Range("E3").Value = "A2"
$ws.Range("F3").Value = "M2"
$ws.Range("G3").Value = "M2"
$ws.Range("J3").Value = "A2"
$ws.Range("K3").Value = "DO"
$ws.Range("L3").Value = "M2"
$ws.Range("M3").Value = "M1"
$ws.Range("O3").Value = "A1"
$ws.Range("S3").Value = "M3"
$ws.Range("T3").Value = "A1"
$ws.Range("U3").Value = "M1"
$ws.Range("V3").Value = "A1"
$ws.Range("W3").Value = "DO"
$ws.Range("Z3").Value = "M2"
$ws.Range("AA3").Value = "M2"
$ws.Range("AB3").Value = "M2"
$ws.Range("AC3").Value = "A2"
$ws.Range("D4").Value = "A1"
$ws.Range("I4").Value = "M1"
$ws.Range("L4").Value = "DO"
$ws.Range("M4").Value = "M1"
$ws.Range("O4").Value = "A1"
$ws.Range("U4").Value = "A1"
$ws.Range("W4").Value = "DO"
$ws.Range("Z4").Value = "M3"
$ws.Range("AA4").Value = "A1"
$ws.Range("AB4").Value = "A1"
$ws.Range("AC4").Value = "M1"
$ws.Range("C5").Value = "A1"
$ws.Range("E5").Value = "A2"
$ws.Range("F5").Value = "A2"
$ws.Range("G5").Value = "M1"
$ws.Range("H5").Value = "DO"
$ws.Range("I5").Value = "A2"
$ws.Range("K5").Value = "DO"
$ws.Range("L5").Value = "M2"
$ws.Range("M5").Value = "M1"
$ws.Range("O5").Value = "M2"
$ws.Range("U5").Value = "M2"
$ws.Range("W5").Value = "M2"
$ws.Range("Z5").Value = "A2"
$ws.Range("AA5").Value = "DO"
$ws.Range("AC5").Value = "M2"
$ws.Range("C6").Value = "M2"
$ws.Range("D6").Value = "A2"
$ws.Range("E6").Value = "M2"
$ws.Range("F6").Value = "M1"
$ws.Range("G6").Value = "M2"
$ws.Range("H6").Value = "A1"
$ws.Range("I6").Value = "M2"
$ws.Range("J6").Value = "A2"
$ws.Range("K6").Value = "M2"
$ws.Range("L6").Value = "DO"
$ws.Range("M6").Value = "M1"
$ws.Range("N6").Value = "M2"
$ws.Range("R6").Value = "M1"
$ws.Range("S6").Value = "A1"
$ws.Range("T6").Value = "DO"
$ws.Range("V6").Value = "M2"
$ws.Range("W6").Value = "M2"
$ws.Range("Z6").Value = "DO"
$ws.Range("AB6").Value = "A1"
$ws.Range("D7").Value = "DO"
$ws.Range("F7").Value = "A1"
$ws.Range("H7").Value = "M3"
$ws.Range("J7").Value = "M3"
$ws.Range("K7").Value = "A1"
$ws.Range("L7").Value = "M1"
$ws.Range("N7").Value = "A1"
$ws.Range("R7").Value = "A1"
$ws.Range("S7").Value = "DO"
$ws.Range("T7").Value = "M3"
$ws.Range("Y7").Value = "M3"
$ws.Range("Z7").Value = "A1"
$ws.Range("AB7").Value = "M1"
$ws.Range("C8").Value = "A2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = "A1"
$ws.Range("F8").Value = "A1"
$ws.Range("G8").Value = "M2"
$ws.Range("H8").Value = "DO"
$ws.Range("I8").Value = "A1"
$ws.Range("J8").Value = "A2"
$ws.Range("K8").Value = "M2"
$ws.Range("L8").Value = "A2"
$ws.Range("M8").Value = "DO"
$ws.Range("O8").Value = "A2"
$ws.Range("Q8").Value = "M2"
$ws.Range("R8").Value = "A2"
$ws.Range("S8").Value = "A1"
$ws.Range("T8").Value = "DO"
$ws.Range("U8").Value = "M1"
$ws.Range("V8").Value = "A2"
$ws.Range("W8").Value = "M1"
$ws.Range("X8").Value = "M1"
$ws.Range("Y8").Value = "A1"
$ws.Range("Z8").Value = "A1"
$ws.Range("AA8").Value = "M1"
$ws.Range("AB8").Value = "M3"
$ws.Range("AC8").Value = "DO"
$ws.Range("B9").Value = "DO"
$ws.Range("C9").Value = "M2"
$ws.Range("D9").Value = "M1"
$ws.Range("F9").Value = "M2"
$ws.Range("G9").Value = "A2"
$ws.Range("H9").Value = "A1"
$ws.Range("J9").Value = "M2"
$ws.Range("K9").Value = "A1"
$ws.Range("L9").Value = "M2"
$ws.Range("N9").Value = "A1"
$ws.Range("Q9").Value = "A2"
$ws.Range("R9").Value = "M1"
$ws.Range("S9").Value = "M1"
$ws.Range("U9").Value = "A2"
$ws.Range("W9").Value = "M2"
$ws.Range("X9").Value = "A2"
$ws.Range("Y9").Value = "A1"
$ws.Range("Z9").Value = "DO"
$ws.Range("AB9").Value = "M2"
$ws.Range("AC9").Value = "A1"
$ws.Range("B10").Value = "A2"
$ws.Range("D10").Value = "A1"
$ws.Range("E10").Value = "A1"
$ws.Range("F10").Value = "DO"
$ws.Range("G10").Value = "M2"
$ws.Range("H10").Value = "M2"
$ws.Range("I10").Value = "A2"
$ws.Range("J10").Value = "M1"
$ws.Range("M10").Value = "DO"
$ws.Range("N10").Value = "M2"
$ws.Range("O10").Value = "M2"
$ws.Range("U10").Value = "M2"
$ws.Range("V10").Value = "DO"
$ws.Range("AB10").Value = "A2"
